$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two of the original numbered rows (C=4 and C=12) are dropped from the list;
# deleting their rows shifts everything below up by two and keeps the
# surviving members of the shared "=D*$G$2" formula block (now E6:E10) intact.
$ws.Rows("11").Delete()
$ws.Rows("7").Delete()

# --- Column A: new running "#" numbering for every main line item ---
$ws.Range("A5").Value = 1
$ws.Range("A11").Value = 2

# --- Row 12: 94KEY OEM PROFILE DOLCH PBT KEYSET (previously row 13) ---
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "94KEY OEM PROFILE DOLCH PBT KEYSET"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = 35

# --- Row 13: ENJOYPBT CMYW/RGBY KEYCAPS (previously row 14) ---
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "ENJOYPBT CMYW/RGBY KEYCAPS"
$ws.Range("C13").Value = "CMYW Blank 24"
$ws.Range("D13").Value = 38.9

# --- Row 14: new "Shiping" row (no Type/column C) ---
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "Shiping"
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = 15
$ws.Range("E14").Formula = "=D14*`$G`$2"

# --- Row 15: POM Laser Engraved 107keys (previously row 12) ---
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = "POM Laser Engraved 107keys"
$ws.Range("C15").Value = "DHL"
$ws.Range("D15").Value = 49.9
$ws.Range("E15").Formula = "=D15*`$G`$2"

# --- Row 17: new summary row ---
$ws.Range("C17").Value = "Sum 1-5"
$ws.Range("E17").Formula = "=SUM(E5:E14)"

# Force a full recalculation so every formula (including the ones whose
# inputs moved around via the row deletes above) carries a fresh cached value.
$excel.CalculateFull()

# --- View state: scroll/selection position matches the edited file ---
$ws.Range("E18").Select()
